# Add a new "2020" data column (O) to the 4.2.2 participation-rate table,
# mirroring the formatting already used by the preceding "2019" column (N).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New year header in O4, copying the look of N4 (bold, right aligned, bordered).
$ws.Range("N4").Copy($ws.Range("O4"))
$ws.Range("O4").Value = 2020

# New data value in O5, copying the look of N5 (right aligned, 1-decimal number format).
$ws.Range("N5").Copy($ws.Range("O5"))
$ws.Range("O5").Value = 83.3

# Move/record the active selection as it ended up after the edit.
$ws.Range("O12").Select()
